{"js": "// Rewrite the six \"Impact\" bullet paragraphs under \"KEY ACHIEVEMENTS AND\n// IMPACT\" into four punchier, impact-focused accomplishment bullets, per the\n// commit: \"Fix Key Achievements to use proper accomplishment statements\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the six old bullet paragraphs by their exact existing text so the\n// script is resilient to any other incidental paragraph-count differences.\nconst oldBullets = [\n  \"\\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n  \"\\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n  \"\\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n  \"\\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\"\n];\n\nconst newBullets = [\n  \"\\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\\u2022 $4.7M savings enabled nonprofit access\",\n  \"\\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"\\u2022 178% accuracy improvement in racial classification algorithms\"\n];\n\nconst items = paragraphs.items;\nlet startIdx = -1;\nfor (let i = 0; i <= items.length - oldBullets.length; i++) {\n  let matched = true;\n  for (let j = 0; j < oldBullets.length; j++) {\n    if (items[i + j].text !== oldBullets[j]) {\n      matched = false;\n      break;\n    }\n  }\n  if (matched) {\n    startIdx = i;\n    break;\n  }\n}\n\nif (startIdx === -1) {\n  throw new Error(\"Could not locate the Key Achievements bullet block to rewrite.\");\n}\n\n// First four old paragraphs become the four new bullets (simple text swap).\nfor (let j = 0; j < newBullets.length; j++) {\n  items[startIdx + j].insertText(newBullets[j], \"Replace\");\n}\n\n// The remaining two old paragraphs (indexes 4 and 5 of the old block) are no\n// longer needed now that we only have four bullets.\nitems[startIdx + 4].delete();\nitems[startIdx + 5].delete();\n\nawait context.sync();\n", "ps1": "# Rewrite the six \"Impact\" bullet paragraphs under \"KEY ACHIEVEMENTS AND\n# IMPACT\" into four punchier, impact-focused accomplishment bullets, per the\n# commit: \"Fix Key Achievements to use proper accomplishment statements\".\n\n$d = $word.ActiveDocument\n\n$oldBullets = @(\n  [char]0x2022 + \" Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\",\n  [char]0x2022 + \" Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n  [char]0x2022 + \" Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  [char]0x2022 + \" Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\",\n  [char]0x2022 + \" Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  [char]0x2022 + \" Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\"\n)\n\n$newBullets = @(\n  [char]0x2022 + \" Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  [char]0x2022 + \" `$4.7M savings enabled nonprofit access\",\n  [char]0x2022 + \" Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  [char]0x2022 + \" 178% accuracy improvement in racial classification algorithms\"\n)\n\n# Locate the six-paragraph block by matching the whole consecutive sequence\n# of old bullet text (the first bullet's text alone is not unique -- it is\n# also used verbatim earlier, under the \"Partner - Siege Analytics\" entry).\n$count = $d.Paragraphs.Count\n$startIndex = -1\nfor ($i = 1; $i -le ($count - $oldBullets.Count + 1); $i++) {\n  $allMatch = $true\n  for ($j = 0; $j -lt $oldBullets.Count; $j++) {\n    $candidate = $d.Paragraphs.Item($i + $j).Range.Text\n    $candidate = $candidate.TrimEnd([char]13, [char]7)\n    if ($candidate -ne $oldBullets[$j]) {\n      $allMatch = $false\n      break\n    }\n  }\n  if ($allMatch) {\n    $startIndex = $i\n    break\n  }\n}\n\nif ($startIndex -eq -1) {\n  throw \"Could not locate the Key Achievements bullet block to rewrite.\"\n}\n\n# First four old paragraphs become the four new bullets (simple text swap).\nfor ($j = 0; $j -lt $newBullets.Count; $j++) {\n  $para = $d.Paragraphs.Item($startIndex + $j)\n  $para.Range.Text = $newBullets[$j]\n}\n\n# The remaining two old paragraphs (the 5th and 6th of the old six-bullet\n# block) are no longer needed now that we only have four bullets. Delete the\n# one at $startIndex+4 twice since paragraphs shift up after each delete.\n$d.Paragraphs.Item($startIndex + 4).Range.Delete()\n$d.Paragraphs.Item($startIndex + 4).Range.Delete()\n"}
